$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "NIG(0.9409090464995247, 0.7041175376194337, 1.0514166398204328, 3.087744930529382)"
$ws.Range("C2").Value = "JSU(-0.8514023414994896, 1.0613851577849713, 4.490502961028408, 4.409352958788498)"
$ws.Range("D2").Value = "NIG(0.5254538891429947, 0.28116030297283046, 2.8132743060584806, 2.8005629058449273)"
$ws.Range("E2").Value = "JSU(-1.084930290429198, 1.1568317259890988, 3.1027869595090563, 5.498227314630439)"
